# Auto-generated edit script applying the Siren_Profits market-data refresh diff.
# Updates currentAveragePrice* / Leve price & profit columns (H-N) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the scheduled runner snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3492.4614
$ws.Range("J18").Value = 2631.5
$ws.Range("L18").Value = 2631.5
$ws.Range("N18").Value = -3199.5
$ws.Range("H40").Value = 3725.3333
$ws.Range("I40").Value = 4338
$ws.Range("J40").Value = 3419
$ws.Range("K40").Value = 4338
$ws.Range("L40").Value = 3419
$ws.Range("M40").Value = -4163
$ws.Range("N40").Value = -3769
$ws.Range("H70").Value = 5426
$ws.Range("J70").Value = 3196.4
$ws.Range("L70").Value = 9589.200000000001
$ws.Range("N70").Value = -10129.2
$ws.Range("H73").Value = 5426
$ws.Range("J73").Value = 3196.4
$ws.Range("L73").Value = 9589.200000000001
$ws.Range("N73").Value = -11461.2
$ws.Range("H94").Value = 50129668
$ws.Range("I94").Value = 71438100
$ws.Range("J94").Value = 410002
$ws.Range("K94").Value = 71438100
$ws.Range("L94").Value = 410002
$ws.Range("M94").Value = -71437649
$ws.Range("N94").Value = -410904
$ws.Range("H107").Value = 7358.2383
$ws.Range("I107").Value = 7056.8335
$ws.Range("J107").Value = 9166.666999999999
$ws.Range("K107").Value = 7056.8335
$ws.Range("L107").Value = 9166.666999999999
$ws.Range("M107").Value = -5136.8335
$ws.Range("N107").Value = -13006.667
$ws.Range("H125").Value = 4568.9165
$ws.Range("J125").Value = 4815.968
$ws.Range("L125").Value = 43343.712
$ws.Range("N125").Value = -48263.712
$ws.Range("H132").Value = 2937.8086
$ws.Range("I132").Value = 2628.4443
$ws.Range("K132").Value = 7885.3329
$ws.Range("M132").Value = -5355.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 7079.9
$ws.Range("I2").Value = 7989.6875
$ws.Range("K2").Value = 7989.6875
$ws.Range("M2").Value = -7876.6875
$ws.Range("H32").Value = 6600.316
$ws.Range("I32").Value = 6866.057
$ws.Range("K32").Value = 6866.057
$ws.Range("M32").Value = -6579.057
$ws.Range("H45").Value = 6913.6924
$ws.Range("I45").Value = 5484.875
$ws.Range("K45").Value = 5484.875
$ws.Range("M45").Value = -5107.875
$ws.Range("H61").Value = 5960.522
$ws.Range("I61").Value = 5808.1
$ws.Range("K61").Value = 5808.1
$ws.Range("M61").Value = -5596.1
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H103").Value = 171666.67
$ws.Range("J103").Value = 171666.67
$ws.Range("L103").Value = 171666.67
$ws.Range("N103").Value = -174010.67
$ws.Range("H109").Value = 60188.5
$ws.Range("J109").Value = 60188.5
$ws.Range("L109").Value = 60188.5
$ws.Range("N109").Value = -62962.5
$ws.Range("H116").Value = 7079.9
$ws.Range("I116").Value = 7989.6875
$ws.Range("K116").Value = 7989.6875
$ws.Range("M116").Value = -5695.6875
$ws.Range("H132").Value = 4133.037
$ws.Range("I132").Value = 3115.5334
$ws.Range("K132").Value = 9346.600199999999
$ws.Range("M132").Value = -6816.600199999999
$ws.Range("H136").Value = 5960.522
$ws.Range("I136").Value = 5808.1
$ws.Range("K136").Value = 17424.3
$ws.Range("M136").Value = -14874.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7079.9
$ws.Range("I3").Value = 7989.6875
$ws.Range("K3").Value = 7989.6875
$ws.Range("M3").Value = -7875.6875
$ws.Range("H86").Value = 4305.7856
$ws.Range("I86").Value = 5269.4707
$ws.Range("J86").Value = 2816.4546
$ws.Range("K86").Value = 5269.4707
$ws.Range("L86").Value = 2816.4546
$ws.Range("M86").Value = -4146.4707
$ws.Range("N86").Value = -5062.4546
$ws.Range("H89").Value = 4305.7856
$ws.Range("I89").Value = 5269.4707
$ws.Range("J89").Value = 2816.4546
$ws.Range("K89").Value = 26347.3535
$ws.Range("L89").Value = 14082.273
$ws.Range("M89").Value = -20731.3535
$ws.Range("N89").Value = -25314.273
$ws.Range("H107").Value = 3989.2188
$ws.Range("I107").Value = 4666.48
$ws.Range("J107").Value = 1570.4286
$ws.Range("K107").Value = 4666.48
$ws.Range("L107").Value = 1570.4286
$ws.Range("M107").Value = -2746.48
$ws.Range("N107").Value = -5410.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10109.827
$ws.Range("I86").Value = 8533.143
$ws.Range("K86").Value = 8533.143
$ws.Range("M86").Value = -7410.143
$ws.Range("H89").Value = 10109.827
$ws.Range("I89").Value = 8533.143
$ws.Range("K89").Value = 42665.715
$ws.Range("M89").Value = -37049.715
$ws.Range("H107").Value = 7231.676
$ws.Range("I107").Value = 9180.406999999999
$ws.Range("K107").Value = 9180.406999999999
$ws.Range("M107").Value = -7260.406999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 11000
$ws.Range("J102").Value = 11000
$ws.Range("L102").Value = 33000
$ws.Range("N102").Value = -37868
$ws.Range("H131").Value = 66667920
$ws.Range("J131").Value = 1799
$ws.Range("L131").Value = 5397
$ws.Range("N131").Value = -15477
$ws.Range("H140").Value = 11771.588
$ws.Range("I140").Value = 12288.5625
$ws.Range("K140").Value = 36865.6875
$ws.Range("M140").Value = -31685.6875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 17000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 17000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H70").Value = 4984.3
$ws.Range("I70").Value = 4499.5
$ws.Range("J70").Value = 5711.5
$ws.Range("K70").Value = 4499.5
$ws.Range("L70").Value = 5711.5
$ws.Range("M70").Value = -4229.5
$ws.Range("N70").Value = -6251.5
$ws.Range("H73").Value = 4984.3
$ws.Range("I73").Value = 4499.5
$ws.Range("J73").Value = 5711.5
$ws.Range("K73").Value = 4499.5
$ws.Range("L73").Value = 5711.5
$ws.Range("M73").Value = -3563.5
$ws.Range("N73").Value = -7583.5
$ws.Range("H107").Value = 678.7368
$ws.Range("I107").Value = 652.7059
$ws.Range("K107").Value = 652.7059
$ws.Range("M107").Value = 1267.2941
$ws.Range("H132").Value = 3699.96
$ws.Range("I132").Value = 3020.7917
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 9062.375100000001
$ws.Range("L132").Value = 60000
$ws.Range("M132").Value = -6532.375100000001
$ws.Range("N132").Value = -65060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4172.8335
$ws.Range("I61").Value = 3096.3704
$ws.Range("K61").Value = 3096.3704
$ws.Range("M61").Value = -2894.3704
$ws.Range("H113").Value = 4172.8335
$ws.Range("I113").Value = 3096.3704
$ws.Range("K113").Value = 3096.3704
$ws.Range("M113").Value = -926.3703999999998
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("I132").Value = 679926.4
$ws.Range("K132").Value = 2039779.2
$ws.Range("M132").Value = -2037249.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 100000
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H112").Value = 40196.75
$ws.Range("J112").Value = 40196.75
$ws.Range("L112").Value = 40196.75
$ws.Range("N112").Value = -43150.75
$ws.Range("H134").Value = 100000
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
